$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.464.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.31%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.658.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.93%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.44%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.48%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3611"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.59%  "

# Row 8
$ws.Range("E8").Value = "  -3.92%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3249"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.78%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.116"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.89%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06980"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.90%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.55%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.886"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.93%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.95%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.652.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.554"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.12%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001042"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.86%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06561"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "76.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.83%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.921"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.64%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.32%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.40%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.420.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.35%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.470"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.308"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -16.95%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.39%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.25%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.842.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.80%  "

# Row 30
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.12%  "

# Row 31
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.169"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.51%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.984"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.46%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.613"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -17.70%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.715"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.72%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08378"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.32%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.96%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.180"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.11%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06034"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.62%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02193"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.53%  "

# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.203"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.67%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2046"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.35%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.181"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.21%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5883"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.10%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.738"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.16%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.03%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5584"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.09%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.16%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.930"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.28%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06895"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.55%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.87%  "
